# TOD-E norms run: split the old "7.0-9.3" raw->ss lookup tab into four
# narrower age-band tabs (7.0-7.5, 7.6-7.11, 8.0-8.5, 8.6-9.3), rescaling
# the "ss" column for the renamed first band and adding three brand-new
# lookup tabs for the remaining bands.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rescale the existing "7.0-9.3" sheet's ss column, then rename it to
#    "7.0-7.5" (it becomes the first of the four split-out age bands).
# ---------------------------------------------------------------------
$wsOld = $wb.Worksheets.Item("7.0-9.3")

$band1 = @(50,53,57,59,61,62,64,65,67,68,70,71,72,74,75,76,78,79,80,82,83,85,86,88,89,91,93,95,98,101,105,111,126)
for ($i = 0; $i -lt $band1.Length; $i++) {
    $wsOld.Cells.Item($i + 2, 2).Value = $band1[$i]
}

$wsOld.Name = "7.0-7.5"

# ---------------------------------------------------------------------
# 2) Build the three new lookup tabs, each shaped like the original:
#    header row (raw / ss, bold + centered via the existing style),
#    then raw = 1..33 in column A with the rescaled ss values in column B.
# ---------------------------------------------------------------------
function Add-NormSheet {
    param(
        [string]$SheetName,
        [int[]]$SsValues,
        [object]$AfterSheet
    )

    $newSheet = $wb.Worksheets.Add($null, $AfterSheet)
    $newSheet.Name = $SheetName

    $newSheet.Range("A1").Value = "raw"
    $newSheet.Range("B1").Value = "ss"
    $newSheet.Range("A1:B1").Font.Bold = $true
    $newSheet.Range("A1:B1").HorizontalAlignment = -4108

    for ($i = 0; $i -lt $SsValues.Length; $i++) {
        $newSheet.Cells.Item($i + 2, 1).Value = $i + 1
        $newSheet.Cells.Item($i + 2, 2).Value = $SsValues[$i]
    }
}

$band2 = @(47,50,53,55,57,58,60,61,63,64,65,66,68,69,70,71,73,74,75,76,78,79,81,82,84,85,87,89,92,94,98,104,117)
$band3 = @(44,47,50,52,53,55,56,58,59,60,61,62,64,65,66,67,68,69,71,72,73,74,76,77,79,80,82,84,86,89,92,97,110)
$band4 = @(41,43,46,48,49,51,52,53,54,56,57,58,59,60,61,62,63,64,65,66,68,69,70,71,73,74,76,78,80,82,85,90,102)

$wsBand1 = $wb.Worksheets.Item("7.0-7.5")
Add-NormSheet "7.6-7.11" $band2 $wsBand1

$wsBand2 = $wb.Worksheets.Item("7.6-7.11")
Add-NormSheet "8.0-8.5" $band3 $wsBand2

$wsBand3 = $wb.Worksheets.Item("8.0-8.5")
Add-NormSheet "8.6-9.3" $band4 $wsBand3

# Restore the originally-active first tab (adding sheets shifts focus to
# the most-recently-created one).
$wb.Worksheets.Item(1).Activate()
